$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3067642
$ws.Range("I132").Value = 6538
$ws.Range("J132").Value = 7003347
$ws.Range("K132").Value = 19614
$ws.Range("L132").Value = 21010041
$ws.Range("M132").Value = -17084
$ws.Range("N132").Value = -21015101

$ws.Range("H138").Value = 1756098.2
$ws.Range("I138").Value = 904.9761999999999
$ws.Range("J138").Value = 3147006.2
$ws.Range("K138").Value = 2714.9286
$ws.Range("L138").Value = 9441018.600000001
$ws.Range("M138").Value = 2425.0714
$ws.Range("N138").Value = -9451298.600000001

$ws.Range("H141").Value = 1474.7028
$ws.Range("I141").Value = 1474.7028
$ws.Range("K141").Value = 4424.1084
$ws.Range("M141").Value = 755.8915999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2142.4443
$ws.Range("I2").Value = 1162.75
$ws.Range("J2").Value = 2926.2
$ws.Range("K2").Value = 1162.75
$ws.Range("L2").Value = 2926.2
$ws.Range("M2").Value = -1049.75
$ws.Range("N2").Value = -3152.2

$ws.Range("H45").Value = 1533.3334
$ws.Range("I45").Value = 1127.2727
$ws.Range("K45").Value = 1127.2727
$ws.Range("M45").Value = -750.2727

$ws.Range("H61").Value = 18219124
$ws.Range("I61").Value = 20429502
$ws.Range("J61").Value = 167702.33
$ws.Range("K61").Value = 20429502
$ws.Range("L61").Value = 167702.33
$ws.Range("M61").Value = -20429290
$ws.Range("N61").Value = -168126.33

$ws.Range("H97").Value = 2404620.8
$ws.Range("I97").Value = 3125669.5
$ws.Range("J97").Value = 1125
$ws.Range("K97").Value = 3125669.5
$ws.Range("L97").Value = 1125
$ws.Range("M97").Value = -3125173.5
$ws.Range("N97").Value = -2117

$ws.Range("H116").Value = 2142.4443
$ws.Range("I116").Value = 1162.75
$ws.Range("J116").Value = 2926.2
$ws.Range("K116").Value = 1162.75
$ws.Range("L116").Value = 2926.2
$ws.Range("M116").Value = 1131.25
$ws.Range("N116").Value = -7514.2

$ws.Range("H136").Value = 18219124
$ws.Range("I136").Value = 20429502
$ws.Range("J136").Value = 167702.33
$ws.Range("K136").Value = 61288506
$ws.Range("L136").Value = 503106.99
$ws.Range("M136").Value = -61285956
$ws.Range("N136").Value = -508206.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2142.4443
$ws.Range("I3").Value = 1162.75
$ws.Range("J3").Value = 2926.2
$ws.Range("K3").Value = 1162.75
$ws.Range("L3").Value = 2926.2
$ws.Range("M3").Value = -1048.75
$ws.Range("N3").Value = -3154.2

$ws.Range("H105").Value = 31251818
$ws.Range("I105").Value = 71430100
$ws.Range("J105").Value = 2043.3334
$ws.Range("K105").Value = 71430100
$ws.Range("L105").Value = 2043.3334
$ws.Range("M105").Value = -71428353
$ws.Range("N105").Value = -5537.3334

$ws.Range("H107").Value = 1541.5385
$ws.Range("I107").Value = 1258.375
$ws.Range("J107").Value = 1994.6
$ws.Range("K107").Value = 1258.375
$ws.Range("L107").Value = 1994.6
$ws.Range("M107").Value = 661.625
$ws.Range("N107").Value = -5834.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1982.8334
$ws.Range("I10").Value = 299
$ws.Range("J10").Value = 3666.6667
$ws.Range("K10").Value = 299
$ws.Range("L10").Value = 3666.6667
$ws.Range("M10").Value = -160
$ws.Range("N10").Value = -3944.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2395.6
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 2482.9167
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 7448.750100000001
$ws.Range("M4").Value = -788
$ws.Range("N4").Value = -7672.750100000001

$ws.Range("H34").Value = 1389.5555
$ws.Range("I34").Value = 916.6667
$ws.Range("J34").Value = 2335.3333
$ws.Range("K34").Value = 2750.0001
$ws.Range("L34").Value = 7005.999899999999
$ws.Range("M34").Value = -2666.0001
$ws.Range("N34").Value = -7173.999899999999

$ws.Range("H39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 302
$ws.Range("I55").Value = 302
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 906
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -729
$ws.Range("N55").ClearContents()

$ws.Range("H131").Value = 1241.6857
$ws.Range("J131").Value = 1263.5
$ws.Range("L131").Value = 3790.5
$ws.Range("N131").Value = -13870.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3681.1155
$ws.Range("I80").Value = 3015.5715
$ws.Range("J80").Value = 3926.3157
$ws.Range("K80").Value = 3015.5715
$ws.Range("L80").Value = 3926.3157
$ws.Range("M80").Value = -2017.5715
$ws.Range("N80").Value = -5922.3157

$ws.Range("H83").Value = 3681.1155
$ws.Range("I83").Value = 3015.5715
$ws.Range("J83").Value = 3926.3157
$ws.Range("K83").Value = 15077.8575
$ws.Range("L83").Value = 19631.5785
$ws.Range("M83").Value = -10085.8575
$ws.Range("N83").Value = -29615.5785

$ws.Range("H123").Value = 27125.47
$ws.Range("J123").Value = 27125.47
$ws.Range("L123").Value = 27125.47
$ws.Range("N123").Value = -32025.47

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 38228.6
$ws.Range("I81").Value = 31000
$ws.Range("J81").Value = 40035.75
$ws.Range("K81").Value = 31000
$ws.Range("L81").Value = 40035.75
$ws.Range("M81").Value = -30002
$ws.Range("N81").Value = -42031.75

$ws.Range("H84").Value = 38228.6
$ws.Range("I84").Value = 31000
$ws.Range("J84").Value = 40035.75
$ws.Range("K84").Value = 93000
$ws.Range("L84").Value = 120107.25
$ws.Range("M84").Value = -88008
$ws.Range("N84").Value = -130091.25

$ws.Range("H93").Value = 1092.7142
$ws.Range("I93").Value = 1018.0909
$ws.Range("K93").Value = 1018.0909
$ws.Range("M93").Value = 229.9091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1261.4546
$ws.Range("I126").Value = 1308.4445
$ws.Range("J126").Value = 1050
$ws.Range("K126").Value = 3925.3335
$ws.Range("L126").Value = 3150
$ws.Range("M126").Value = -1455.3335
$ws.Range("N126").Value = -8090

$ws.Range("H132").Value = 34251.164
$ws.Range("I132").Value = 27188.986
$ws.Range("J132").Value = 65407.824
$ws.Range("K132").Value = 81566.958
$ws.Range("L132").Value = 196223.472
$ws.Range("M132").Value = -79036.958
$ws.Range("N132").Value = -201283.472
